$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-08-18 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-19 Monday", 2)

# Update each division expression in the practice table.
# Cells are addressed directly by (row, column) to avoid ambiguity from
# duplicate / overlapping expression text between old and new values.
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="35÷8="},
    @{Row=1;  Col=2; Text="73÷4="},
    @{Row=1;  Col=3; Text="67÷2="},
    @{Row=1;  Col=4; Text="57÷4="},
    @{Row=1;  Col=5; Text="28÷5="},

    @{Row=5;  Col=1; Text="50÷6="},
    @{Row=5;  Col=2; Text="39÷3="},
    @{Row=5;  Col=3; Text="22÷6="},
    @{Row=5;  Col=4; Text="37÷3="},
    @{Row=5;  Col=5; Text="19÷3="},

    @{Row=9;  Col=1; Text="83÷7="},
    @{Row=9;  Col=2; Text="37÷7="},
    @{Row=9;  Col=3; Text="92÷6="},
    @{Row=9;  Col=4; Text="86÷7="},
    @{Row=9;  Col=5; Text="70÷9="},

    @{Row=13; Col=1; Text="85÷3="},
    @{Row=13; Col=2; Text="67÷2="},
    @{Row=13; Col=3; Text="93÷2="},
    @{Row=13; Col=4; Text="11÷7="},
    @{Row=13; Col=5; Text="80÷3="},

    @{Row=17; Col=1; Text="90÷3="},
    @{Row=17; Col=2; Text="66÷4="},
    @{Row=17; Col=3; Text="66÷9="},
    @{Row=17; Col=4; Text="49÷2="},
    @{Row=17; Col=5; Text="50÷5="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
